$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.695.82'
$ws.Range("E2").Value = '  -0.77%  '

$ws.Range("D3").Value = '2.227.62'
$ws.Range("E3").Value = '  -2.08%  '

$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.63'
$ws.Range("E5").Value = '  -1.87%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.94'
$ws.Range("E6").Value = '  -5.23%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.568'
$ws.Range("E7").Value = '  -3.46%  '

$ws.Range("E8").Value = '  +0.19%  '

$ws.Range("E9").Value = '  -7.04%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.66'
$ws.Range("E10").Value = '  -7.51%  '

$ws.Range("E11").Value = '  -2.33%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.33'
$ws.Range("E12").Value = '  -6.79%  '

$ws.Range("E13").Value = '  -3.12%  '

$ws.Range("D14").Value = '2.566.79'
$ws.Range("E14").Value = '  -2.14%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '2.234.25'
$ws.Range("E15").Value = '  -2.16%  '

$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.835'
$ws.Range("E16").Value = '  -4.42%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.98'
$ws.Range("E17").Value = '  -3.74%  '

$ws.Range("D18").Value = '43.600.77'
$ws.Range("E18").Value = '  -0.78%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.99'
$ws.Range("E19").Value = '  -9.40%  '

$ws.Range("D20").Value = '0.0₃0961'
$ws.Range("E20").Value = '  -3.58%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.28'
$ws.Range("E21").Value = '  -5.70%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.01'
$ws.Range("E22").Value = '  -1.72%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.83'
$ws.Range("E23").Value = '  -1.27%  '

$ws.Range("E24").Value = '  -7.98%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.01'
$ws.Range("E25").Value = '  -8.69%  '

$ws.Range("E26").Value = '  +0.15%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.94'
$ws.Range("E27").Value = '  -3.14%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.18'
$ws.Range("E28").Value = '  -1.95%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.37'
$ws.Range("E29").Value = '  -8.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.96'
$ws.Range("E30").Value = '  -8.44%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '156.99'
$ws.Range("E31").Value = '  -3.21%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.82'
$ws.Range("E32").Value = '  -3.27%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0824'
$ws.Range("E33").Value = '  -6.20%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.64'
$ws.Range("E34").Value = '  -2.76%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.11'
$ws.Range("E35").Value = '  -5.39%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.108'
$ws.Range("E36").Value = '  +0.00%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.87'
$ws.Range("E37").Value = '  -7.99%  '

$ws.Range("E38").Value = '  -3.58%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.32'
$ws.Range("E39").Value = '  -2.40%  '

$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.51'
$ws.Range("E40").Value = '  -9.04%  '

$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.99'
$ws.Range("E41").Value = '  -11.89%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0305'
$ws.Range("E42").Value = '  -6.47%  '

$ws.Range("E43").Value = '  +0.21%  '

$ws.Range("D44").Value = '1.704.37'
$ws.Range("E44").Value = '  -4.19%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '83.71'
$ws.Range("E45").Value = '  -1.26%  '

$ws.Range("E46").Value = '  -7.10%  '

$ws.Range("E47").Value = '  -5.59%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '101.26'
$ws.Range("E48").Value = '  -3.15%  '

$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '70.95'
$ws.Range("E49").Value = '  -5.08%  '

$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.62'
$ws.Range("E50").Value = '  +0.42%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '55.70'
$ws.Range("E51").Value = '  -6.49%  '
